$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "ghymat" header in B1 to "Ghymatkoly"
$ws.Range("B1").Value = "Ghymatkoly"

# Refresh the weekly totals in column B (rows 2-8)
$ws.Range("B2").Value = 189291715
$ws.Range("B3").Value = 272055428
$ws.Range("B4").Value = 271232758
$ws.Range("B5").Value = 298567450
$ws.Range("B6").Value = 309812168
$ws.Range("B7").Value = 283373610
$ws.Range("B8").Value = 293839298

# New column widths: column A gets an explicit width, column B is narrowed a bit
$ws.Columns.Item(1).ColumnWidth = 12.65
$ws.Columns.Item(2).ColumnWidth = 14.33
